# DeudoresPrueba.xlsx update
#
# The author re-sorted the debtor ledger alphabetically by client name and
# appended new collections transactions (new clients + new rows), bringing
# the table from 28 data rows (A2:E29) to 38 data rows (A2:E39).
#
# Note: the workbook's x15ac:absPath (the folder Excel remembers the file
# was last opened from) is written by Excel itself from the OS save path;
# it is not a property exposed on the Workbook/Application object model,
# so it cannot be set from automation code (real VBA/COM can't touch it
# either) and is left alone here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final table contents, sorted by Cliente (column B), row-for-row as it
# should appear starting at A2. Columns: Consecutivo, Cliente, Fecha
# (date serial), Valor, Pagado (always FALSE in this sheet).
$table = @(
    @(1,  "ALISO",                  45996, 196000),
    @(2,  "CALDAS WOK",              45997, 85000),
    @(3,  "CAMILIN",                 45997, 166000),
    @(4,  "CAMPO VERDE TOCANCIPA",   45995, 635000),
    @(5,  "CAMPO VERDE ZIPAQUIRA",   45988, 64200),
    @(6,  "CAMPO VERDE ZIPAQUIRA",   45995, 620000),
    @(7,  "CARNES JOHANA",           45993, 176000),
    @(8,  "CARNILANDIA",             45996, 82000),
    @(9,  "CARNILANDIA",             45997, 300700),
    @(10, "CARNIVOROS",              45959, 200000),
    @(11, "CIMARRON DORADO",         45992, 407000),
    @(12, "CIMARRON DORADO",         45996, 298700),
    @(13, "COCINA CHINA",            45994, 85000),
    @(14, "DARWIN FUTBOL",           45921, 200000),
    @(15, "DAVIDCITO",               45947, 100000),
    @(16, "EL RES",                  45997, 240000),
    @(17, "EL RUBY",                 45992, 85100),
    @(18, "LA PAMPA",                45994, 249000),
    @(19, "LA SELECTA",              45912, 82000),
    @(20, "LOS PAISANOS",            45995, 438500),
    @(21, "MARIANA",                 45650, 171900),
    @(22, "MERKA FRUVER ALEJANDRO",  45988, 60900),
    @(23, "MERKA FRUVER ALEJANDRO",  45995, 893700),
    @(24, "MERKA FRUVER DEXI",       45988, 115400),
    @(25, "MERKA FRUVER DEXI",       45995, 339000),
    @(26, "NEVADA",                  45996, 229000),
    @(27, "NOVILLON SAN MATEO",      45971, 83000),
    @(28, "PARAÍSO FUNZA",           45996, 202000),
    @(29, "PINILLA",                 45924, 16000),
    @(30, "PINILLA",                 45931, 166000),
    @(31, "PINILLA SOACHA",          45993, 129000),
    @(32, "PLAZA JESSICA",           45993, 621000),
    @(33, "PLAZA JESSICA",           45995, 1580300),
    @(34, "PORTAL ZIPA",             45995, 66400),
    @(35, "SANTANDER MADRID",        45996, 63000),
    @(36, "SANTANDER SUR",           45997, 250700),
    @(37, "SANTANDER SUR",           45993, 80000),
    @(38, "VNZLNO PUNTA ANCA",       45992, 82000)
)

$row = 2
foreach ($rec in $table) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]

    $dateCell = $ws.Cells.Item($row, 3)
    $dateCell.Value = $rec[2]
    # Keep the same custom date format already used by the column
    # (style index 2 in the original file, format "yyyy-mm-dd").
    $dateCell.NumberFormat = "yyyy\-mm\-dd"

    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $false

    $row = $row + 1
}

# Match the author's final cursor position.
$ws.Range("J31").Select() | Out-Null
